$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1; 3 = 0; 4 = 0; 5 = 1; 6 = 0; 7 = 0; 8 = 0; 9 = 1; 10 = 0;
    11 = 1; 12 = 0; 13 = 0; 14 = 1; 15 = 2; 16 = 0; 17 = 2; 18 = 2; 19 = 0;
    20 = 1; 21 = 0; 22 = 2; 23 = 2; 24 = 3; 25 = 1; 26 = 0; 27 = 1; 28 = 0;
    29 = 3; 30 = 1; 31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 0; 36 = 1; 37 = 1;
    38 = 2; 39 = 0; 40 = 2; 41 = 0; 42 = 0; 43 = 2; 44 = 3; 45 = 0; 46 = 2;
    47 = 2; 48 = 3; 49 = 3; 50 = 4; 51 = 2; 52 = 0; 53 = 0; 54 = 1; 55 = 2;
    56 = 2; 57 = 3; 58 = 3; 59 = 3; 60 = 1; 61 = 2; 62 = 0; 63 = 1; 64 = 0;
    65 = 1; 66 = 1; 67 = 1; 68 = 1; 69 = 3; 70 = 1; 71 = 4; 72 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
